$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.322.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.76%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.80"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4523"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.39%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.47"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -8.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.013"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.37"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.862.88"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.905"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.123"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001027"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06556"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.39%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.529"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.311.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.067.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.90"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.065"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.452"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.28%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09300"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9339"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.459"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.602"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.265"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02221"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05987"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.094"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -11.54%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5907"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1884"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.10"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.271"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5620"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.07"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.372"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.915"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06753"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.09%  "
